$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.656.12"
$ws.Range("E2").Value = "  +1.12%  "

$ws.Range("D3").Value = "3.395.93"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.31"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.90"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.62"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.36%  "

$ws.Range("E10").Value = "  -0.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.386"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.21%  "

$ws.Range("D12").Value = "3.975.73"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("E13").Value = "  -0.21%  "

$ws.Range("E14").Value = "  +0.94%  "

$ws.Range("D15").Value = "3.388.24"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("D17").Value = "61.710.30"
$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.14"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.67"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.12"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.02"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.64"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.549"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("E25").Value = "  -3.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.181"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.42%  "

$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.40"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.99"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.15"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.52%  "

$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.39"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.94"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "168.37"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.10"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.35%  "

$ws.Range("D37").Value = "3.429.11"

$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("E39").Value = "  -1.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.12"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.57%  "

$ws.Range("E41").Value = "  +0.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.45"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.40%  "

$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("E44").Value = "  +2.47%  "

$ws.Range("D45").Value = "2.477.78"
$ws.Range("E45").Value = "  +1.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.67"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.66"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.97%  "

$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0264"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("E50").Value = "  -5.15%  "

$ws.Range("E51").Value = "  -1.40%  "

